$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column (D) cells that hold numeric-looking text
# (e.g. "245.33") are kept as plain text, matching the source data,
# instead of being auto-converted into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.961.47"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.038.75"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.33"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("E6").Value = "  -1.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.44"
$ws.Range("E7").Value = "  -1.85%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -2.23%  "

$ws.Range("E10").Value = "  -2.20%  "

$ws.Range("E11").Value = "  +2.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.31"
$ws.Range("E12").Value = "  -5.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.887"
$ws.Range("E13").Value = "  +8.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.336.77"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.64"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.969.72"
$ws.Range("E16").Value = "  -4.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.27"
$ws.Range("E17").Value = "  +1.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.971.71"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.47"
$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0880"
$ws.Range("E20").Value = "  -2.12%  "

$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.12"
$ws.Range("E22").Value = "  -1.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.55"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.57"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("E27").Value = "  -2.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.88"
$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.50"
$ws.Range("E29").Value = "  +14.93%  "

$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  -3.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.75"
$ws.Range("E32").Value = "  +5.22%  "

$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.85"
$ws.Range("E35").Value = "  +6.31%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0855"
$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("E38").Value = "  -4.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.22"
$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("E41").Value = "  -0.42%  "

$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("E43").Value = "  -16.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.81"
$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.82"
$ws.Range("E45").Value = "  -5.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.292.03"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("E47").Value = "  -4.39%  "

$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.67"
$ws.Range("E49").Value = "  +6.54%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.223.74"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.66"
$ws.Range("E51").Value = "  -2.14%  "
